$d = $word.ActiveDocument

# Locate the target paragraph: "Created Partial Views. Created the _Layout and _ViewStart views. Created the _ViewImports file."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Created Partial Views.*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Created Partial Views' paragraph."
}

$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# Paragraph 1 (rewritten): splits the original two sentences into many runs.
$para1 = '<w:p>' + $pPr +
    '<w:r><w:t>Created Partial Views</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> by c</w:t></w:r>' +
    '<w:r><w:t>reat</w:t></w:r>' +
    '<w:r><w:t>ing</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> the _Layout and _ViewStart views</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">creating </w:t></w:r>' +
    '<w:r><w:t>the _ViewImports file.</w:t></w:r>' +
    '</w:p>'

# Paragraph 2 (new): describes installing Bootstrap/JQuery, with proofErr spell-check markers.
$para2 = '<w:p>' + $pPr +
    '<w:r><w:t>Installed</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Bootstrap and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>JQuery</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>to the project</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Added </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>bootstrao</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>css</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>jquery</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> to the _Layout view.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Added Bootstrap and made the list page more appealing.</w:t></w:r>' +
    '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $para1 + $para2 +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML on the paragraph's own range replaces that paragraph (including its
# paragraph mark/pPr) with the supplied XML, so one source paragraph becomes two.
$target.Range.InsertXML($xml)

Write-Output "Paragraph split and Bootstrap/JQuery paragraph inserted."
